# Total Forces fixed - calculated by openmm
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-8), replacing previous single data row
$data = @(
    @("130, 455, 780",          1, "130",       "5269"),
    @("1073, 1105",             1, "1105",      "5131"),
    @("130, 1073, 1105",        1, "1105",      "5399"),
    @("423, 748, 780, 1073",    2, "780, 780",  "5677, 5887"),
    @("423, 1105, 1105",        1, "1105",      "5331"),
    @("423, 1073, 1105",        1, "1105",      "5433"),
    @("98, 130, 748, 1073",     1, "130",       "5582")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
